$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Ideal Time" (column E) estimates for the first user story's tasks ---
$ws.Range("E2:E6").HorizontalAlignment = -4108
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 0.5
$ws.Range("E4").Value = 0.5
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 1

# --- Record Day2 progress for tasks 1.3 / 1.4 / 1.5 (rows 6-8) ---
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 2

# --- Insert a new "6." task row before the existing "Est. Time" totals row ---
$ws.Rows("12:12").Insert()

$ws.Range("E12").HorizontalAlignment = -4108
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1.5

# --- Update the "Est. Time" row's Day2 total to reflect the new work ---
$ws.Range("C13").Value = 11

# --- Add a new user story #6 at the bottom of the sheet ---
$ws.Range("A36").Value = "6. Implementering av inläsnig av textfil för sparning av medlemmar."

# --- Update the active selection / scroll position ---
$ws.Activate()
$ws.Range("H12").Select()
